# Update the "想去人数" (interested-count) column F values across the four
# worksheets to reflect newly scraped numbers (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 8451
$ws1.Range("F7").Value  = 2388
$ws1.Range("F11").Value = 1045
$ws1.Range("F14").Value = 315
$ws1.Range("F15").Value = 2100
$ws1.Range("F21").Value = 1412
$ws1.Range("F22").Value = 633
$ws1.Range("F23").Value = 1660

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 56
$ws2.Range("F37").Value = 52
$ws2.Range("F39").Value = 303

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F16").Value = 2726
$ws3.Range("F18").Value = 638

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 8451
$ws4.Range("F12").Value = 2726
$ws4.Range("F14").Value = 1045
$ws4.Range("F21").Value = 315
$ws4.Range("F28").Value = 1412
$ws4.Range("F32").Value = 633
$ws4.Range("F35").Value = 1660
$ws4.Range("F43").Value = 52
$ws4.Range("F44").Value = 303
